{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the changes described by the commit diff:\n//   1. Update the \"Generated using git revision\" value.\n//   2. In the Methods paragraph, insert a new sentence about right-censoring\n//      and drop the now-redundant \"(ongoing reviews)\" parenthetical.\n//   3. Move the \"We updated the preprint version...\" text out of the Methods\n//      paragraph (shortening it to end at \"...survivor functions.\") and use\n//      it to replace the \"TODO: Describe any protocol deviations.\" sentence\n//      in the Appendix 1 paragraph.\n\nconst body = context.document.body;\n\n// --- 1. Git revision ------------------------------------------------------\nconst revResults = body.search(\"23e6c34\", { matchCase: true, matchWholeWord: false });\nrevResults.load(\"items\");\nawait context.sync();\nfor (const item of revResults.items) {\n  item.insertText(\"6024469\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2. Methods paragraph: censoring sentence ------------------------------\nconst oldCensor =\n  \"All analyses accounted for right-censored outcomes (ongoing reviews) and for nonrandom\";\nconst newCensor =\n  \"Ongoing reviews were right censored at the end of data collection (31 January 2023). \" +\n  \"All analyses accounted for right-censored outcomes and for nonrandom\";\n\nconst censorResults = body.search(oldCensor, { matchCase: true });\ncensorResults.load(\"items\");\nawait context.sync();\nfor (const item of censorResults.items) {\n  item.insertText(newCensor, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 3. Move the protocol-deviation text out of Methods --------------------\nconst movedText =\n  \"We updated the preprint version of the protocol during data extraction but before starting \" +\n  \"the analysis or unblinding the statistician (CJR) to redefine the comparisons in terms of \" +\n  \"under- and overuse of machine learning (TODO: Cite revision). However, only two reviews were \" +\n  \"judged to have under- or overused machine learning, so it was not possible to perform the \" +\n  \"revised analyses. We therefore performed and report the analyses as originally planned.\";\n\nconst oldMethodsTail = \" \" + movedText;\nconst methodsTailResults = body.search(oldMethodsTail, { matchCase: true });\nmethodsTailResults.load(\"items\");\nawait context.sync();\nfor (const item of methodsTailResults.items) {\n  item.insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 4. Replace the Appendix 1 TODO with the moved text ---------------------\nconst todoResults = body.search(\"TODO: Describe any protocol deviations.\", { matchCase: true });\ntodoResults.load(\"items\");\nawait context.sync();\nfor (const item of todoResults.items) {\n  item.insertText(movedText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# Applies the changes described by the commit diff:\n#   1. Update the \"Generated using git revision\" value.\n#   2. In the Methods paragraph, insert a new sentence about right-censoring\n#      and drop the now-redundant \"(ongoing reviews)\" parenthetical.\n#   3. Move the \"We updated the preprint version...\" text out of the Methods\n#      paragraph (shortening it to end at \"...survivor functions.\") and use\n#      it to replace the \"TODO: Describe any protocol deviations.\" sentence\n#      in the Appendix 1 paragraph.\n\n$d = $word.ActiveDocument\n\n# --- 1. Git revision --------------------------------------------------------\n$d.Content.Find.Execute(\"23e6c34\", $false, $false, $false, $false, $false, $true, 1, $false, \"6024469\", 2)\n\n# --- 2. Methods paragraph: censoring sentence -------------------------------\n$oldCensor = \"All analyses accounted for right-censored outcomes (ongoing reviews) and for nonrandom\"\n$newCensor = \"Ongoing reviews were right censored at the end of data collection (31 January 2023). All analyses accounted for right-censored outcomes and for nonrandom\"\n$d.Content.Find.Execute($oldCensor, $false, $false, $false, $false, $false, $true, 1, $false, $newCensor, 2)\n\n# --- 3. Move the protocol-deviation text out of Methods ---------------------\n$movedText = \"We updated the preprint version of the protocol during data extraction but before starting the analysis or unblinding the statistician (CJR) to redefine the comparisons in terms of under- and overuse of machine learning (TODO: Cite revision). However, only two reviews were judged to have under- or overused machine learning, so it was not possible to perform the revised analyses. We therefore performed and report the analyses as originally planned.\"\n\n$oldMethodsTail = \" \" + $movedText\n$d.Content.Find.Execute($oldMethodsTail, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n# --- 4. Replace the Appendix 1 TODO with the moved text ---------------------\n$d.Content.Find.Execute(\"TODO: Describe any protocol deviations.\", $false, $false, $false, $false, $false, $true, 1, $false, $movedText, 2)\n"}
